$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert new column "xi_3" at D, shifting the old "E" header to column E ---
$ws.Range("D1").Value = "xi_3"
$ws.Range("E1").Value = "E"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats: copy D1 formatting (bold, border, centered) to E1
$excel.CutCopyMode = 0

# --- Data values (rows 2..42), columns B (xi_1), C (xi_2), D (xi_3, new), E (E) ---
$colB = @("6.8","2.04","5.372","3.0396","4.67228","3.529404","4.3294172","3.76940796","4.161414428","3.8870099004","4.07909306972","3.944634851196","4.0387556041628","3.97287107708604","4.01899024603977","3.98670682777216","4.00930522055949","3.99348634560836","4.00455955807415","3.9968083093481","4.00223418345633","3.99843607158057","4.0010947498936","3.99923367507448","4.00053642744787","3.99962450078649","4.00026284944945","3.99981600538538","4.00012879623023","3.99990984263884","4.00006311015281","3.99995582289303","4.00003092397488","3.99997835321758","4.00001515274769","3.99998939307662","4.00000742484637","3.99999480260754","4.00000363817472","3.9999974532777","4.00000178270561")
$colC = @("3.4","1.02","2.686","1.5198","2.33614","1.764702","2.1647086","1.88470398","2.080707214","1.9435049502","2.03954653486","1.972317425598","2.0193778020814","1.98643553854302","2.00949512301989","1.99335341388608","2.00465261027974","1.99674317280418","2.00227977903707","1.99840415467405","2.00111709172817","1.99921803579028","2.0005473749468","1.99961683753724","2.00026821372393","1.99981225039325","2.00013142472473","1.99990800269269","2.00006439811512","1.99995492131942","2.00003155507641","1.99997791144652","2.00001546198744","1.99998917660879","2.00000757637385","1.99999469653831","2.00000371242318","1.99999740130377","2.00000181908736","1.99999872663885","2.00000089135281")
$colD = @("1.13333333333333","0.34","0.895333333333333","0.5066","0.778713333333333","0.588234","0.721569533333333","0.62823466","0.693569071333333","0.6478349834","0.679848844953333","0.657439141866","0.673125934027133","0.66214517951434","0.669831707673295","0.664451137962027","0.668217536759915","0.665581057601393","0.667426593012358","0.666134718224683","0.667039030576056","0.666406011930095","0.666849124982267","0.666538945845746","0.666756071241311","0.666604083464416","0.666710474908242","0.666636000897564","0.666688132705039","0.666651640439806","0.666677185025469","0.666659303815505","0.66667182066248","0.666663058869598","0.666669192124615","0.666664898846103","0.666667904141061","0.66666580043459","0.66666727302912","0.666666242212949","0.666666963784269")
$colE = @("6.8","4.76","3.332","2.3324","1.63268","1.142876","0.8000132","0.56000924","0.392006468","0.274404527600001","0.192083169320001","0.134458218524001","0.0941207529668002","0.0658845270767596","0.0461191689537319","0.0322834182676126","0.0225983927873288","0.0158188749511301","0.0110732124657908","0.0077512487260533","0.0054258741082366","0.0037981118757652","0.0026586783130362","0.0018610748191259","0.0013027523733883","0.0009119266613719","0.0006383486629601","0.0004468440640721","0.0003127908448501","0.0002189535913945","0.0001532675139763","0.0001072872597838","7.51010818493825e-05","5.2570757294923e-05","3.67995301067126e-05","2.57596710753205e-05","1.8031769752902e-05","1.26222388270314e-05","8.83556717878875e-06","6.18489702519653e-06","4.32942791794844e-06")

# Store as text (matches the workbook-wide convention of numbers-as-text for this table)
$ws.Range("B2:E42").NumberFormat = "@"

for ($i = 0; $i -lt $colB.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $colB[$i]
    $ws.Cells.Item($r, 3).Value = $colC[$i]
    $ws.Cells.Item($r, 4).Value = $colD[$i]
    $ws.Cells.Item($r, 5).Value = $colE[$i]
}

# --- New row 42: index column A42 = "41" ---
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = "41"
